$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row above the current row 557, shifting rows 557-657
# down to 558-658 (this also extends the sheet dimension to A1:R658).
$ws.Rows.Item(557).EntireRow.Insert()

# Populate the newly inserted row 557 with the new weekly price record.
$ws.Cells.Item(557, 1).Value = 3
$ws.Cells.Item(557, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(557, 3).Value = "Coquimbo"
$ws.Cells.Item(557, 4).Value = 44995
$ws.Cells.Item(557, 5).Value = 5
$ws.Cells.Item(557, 6).Value = 100112037
$ws.Cells.Item(557, 7).Value = "Cebollín"
$ws.Cells.Item(557, 8).Value = "Sin especificar"
$ws.Cells.Item(557, 9).Value = "Primera"
$ws.Cells.Item(557, 10).Value = 230
$ws.Cells.Item(557, 11).Value = 3000
$ws.Cells.Item(557, 12).Value = 3500
$ws.Cells.Item(557, 13).Value = 3239
$ws.Cells.Item(557, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(557, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(557, 16).Value = 90
$ws.Cells.Item(557, 17).Value = 36
$ws.Cells.Item(557, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same date-formatted style as the rest
# of column D.
$ws.Cells.Item(557, 4).NumberFormat = $ws.Cells.Item(556, 4).NumberFormat
